$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift existing rows down by one, preserving the header row's
#     row-level style (s="1" customFormat="1") and its cell styles
#     (A1/B1 were bold+centered, C1/D1/E1 were bold only).
$ws.Rows(1).Insert()

# --- Clear the stray cells left behind by the old row 3 (now row 4)
#     that are not part of the new layout (old Passthrough/bool cells).
$ws.Range("C4").ClearContents()
$ws.Range("D4").ClearContents()

# --- Header row (row 2): Operation / File Path / Row start / Row end /
#     Wave Name / Passthrough / Host Name / Parallel Count / Policy Name /
#     Test Mode / Datastore
$ws.Range("A2").Value = "Operation"
$ws.Range("B2").Value = "File Path Which Includes Data Related To Operation"
$ws.Range("C2").Value = "Row start from excel file"
$ws.Range("D2").Value = "Row end from excel file"
$ws.Range("E2").Value = "Wave Name"
$ws.Range("F2").Value = "Passthrough (True/False)"
$ws.Range("G2").Value = "Host Name"
$ws.Range("H2").Value = "Parallel Count"
$ws.Range("I2").Value = "Policy Name"
$ws.Range("J2").Value = "Test Mode (True/False)"
$ws.Range("K2").Value = "Datastore"

# Row 3
$ws.Range("A3").Value = "Add wave with upload file"
$ws.Range("B3").Value = "C:\Users\Pranav Pawar\PycharmProjects\RMM_DataDriven\TestData\OneForAll\Second Flow.csv"

# Row 4
$ws.Range("A4").Value = "Add vCenter"
$ws.Range("B4").Value = "C:\Users\Pranav Pawar\PycharmProjects\RMM_DataDriven\TestData\OneForAll\addVcenter.xlsx"

# Row 5
$ws.Range("A5").Value = "Set Autoprovision"
$ws.Range("B5").Value = "C:\Users\Pranav Pawar\PycharmProjects\RMM_DataDriven\TestData\OneForAll\setAutoprovisionAndNIC.xlsx"

# Row 6
$ws.Range("A6").Value = "Bulk Edit Sync options for wave"
$ws.Range("B6").Value = "C:\Users\Pranav Pawar\PycharmProjects\RMM_DataDriven\TestData\OneForAll\bulkEditOptions.xlsx"

# Row 7
$ws.Range("A7").Value = "Edit Sync Options"
$ws.Range("B7").Value = "C:\Users\Pranav Pawar\PycharmProjects\RMM_DataDriven\TestData\OneForAll\editSyncOptions.xlsx"

# Row 8
$ws.Range("A8").Value = "Change Datastore for all waves"
$ws.Range("E8").Value = "Second Flow"
$ws.Range("K8").Value = "esx09-datastore2"

# Row 9
$ws.Range("A9").Value = "Bulk Edit Sync options for windows waves"
$ws.Range("B9").Value = "C:\Users\Pranav Pawar\PycharmProjects\RMM_DataDriven\TestData\OneForAll\bulkEditOptionsWindows.xlsx"

# Row 10
$ws.Range("A10").Value = "Set Parallel Count"
$ws.Range("E10").Value = "Second Flow"
$ws.Range("H10").Value = 4

# Row 11
$ws.Range("A11").Value = "Start wave and verify"
$ws.Range("E11").Value = "Second Flow"

# Row 12
$ws.Range("A12").Value = "Check Wave Status"
$ws.Range("E12").Value = "Second Flow"

# --- Header row styling, reusing the two pre-existing cell styles
#     (bold+centered on A2:D2, bold-only on E2:K2) via a format-only
#     copy/paste so no new style entries are created.
$ws.Range("A2").Copy()
$ws.Range("C2:D2").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("F2:K2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column widths (best-effort match of bestFit widths)
$ws.Columns(3).ColumnWidth = 20.666666666666668
$ws.Columns(4).ColumnWidth = 20.0
$ws.Columns(5).ColumnWidth = 10.333333333333334
$ws.Columns(6).ColumnWidth = 21.5
$ws.Columns(7).ColumnWidth = 9.333333333333334
$ws.Columns(8).ColumnWidth = 11.833333333333334
$ws.Columns(9).ColumnWidth = 10.666666666666666
$ws.Columns(10).ColumnWidth = 19.833333333333332

# --- Selection
[void]$ws.Range("B6").Select()
